$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 83030000
$ws.Range("B3").Value = "Kalpana"
$ws.Range("C3").Value = 45819
$ws.Range("C3").NumberFormat = $ws.Range("C2").NumberFormat
$ws.Range("D3").Value = 1000
$ws.Range("E3").Value = 11000
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 6000
$ws.Range("H3").Value = 2000
$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 5000

$ws.Range("E6").Select() | Out-Null
